$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.000.45"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "1.819.59"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  -0.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.03"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4508"
$ws.Range("E7").Value = "  +6.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3696"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07275"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8552"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.73"
$ws.Range("E11").Value = "  -1.21%  "

$ws.Range("D12").Value = "1.807.57"
$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07100"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.24"
$ws.Range("E15").Value = "  +4.71%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.323"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.96"
$ws.Range("E20").Value = "  -0.59%  "

$ws.Range("D21").Value = "26.920.39"
$ws.Range("E21").Value = "  -1.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.164"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.93"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.985"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.59"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.215"
$ws.Range("E26").Value = "  +4.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  +0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.243"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.34"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08866"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.180"
$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7512"
$ws.Range("E32").Value = "  -0.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.963"
$ws.Range("E33").Value = "  +5.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.437"
$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.098"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05239"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5304"
$ws.Range("E39").Value = "  +5.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.174"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.876"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("E42").Value = "  +0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5237"
$ws.Range("E43").Value = "  +10.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.514"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.62"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.975"
$ws.Range("E46").Value = "  +9.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.59"
$ws.Range("E47").Value = "  -1.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9999"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06365"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("E51").Value = "  +0.47%  "
